# Updates the crypto price/volume table on Sheet1 to reflect the latest
# coin-ranking snapshot (values refreshed by the scheduled GitHub Action).
#
# Column D ("Price") holds values such as "30.345.10" or "1.000" that look
# numeric but must stay literal text (they use "." as a thousands/
# formatting separator and rely on trailing zeros). Setting a plain .Value
# on such a string lets Excel auto-convert it to a real number and lose the
# exact text, so each Price cell is briefly switched to the Text ("@")
# number format, assigned, and then restored to its original style so the
# workbook keeps its original (unstyled) look for these cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

# Row 2
Set-TextValue $ws.Cells.Item(2, 4) '30.345.10'
$ws.Cells.Item(2, 5).Value = '  +0.04%  '

# Row 3
Set-TextValue $ws.Cells.Item(3, 4) '1.934.78'
$ws.Cells.Item(3, 5).Value = '  +0.18%  '

# Row 4
Set-TextValue $ws.Cells.Item(4, 4) '1.000'
$ws.Cells.Item(4, 5).Value = '  -0.11%  '

# Row 5
Set-TextValue $ws.Cells.Item(5, 4) '0.7587'
$ws.Cells.Item(5, 5).Value = '  +6.11%  '

# Row 6
Set-TextValue $ws.Cells.Item(6, 4) '244.76'
$ws.Cells.Item(6, 5).Value = '  -2.76%  '

# Row 7
Set-TextValue $ws.Cells.Item(7, 4) '0.9999'
$ws.Cells.Item(7, 5).Value = '  -0.15%  '

# Row 8
Set-TextValue $ws.Cells.Item(8, 4) '27.65'
$ws.Cells.Item(8, 5).Value = '  +0.81%  '

# Row 9
Set-TextValue $ws.Cells.Item(9, 4) '0.3177'
$ws.Cells.Item(9, 5).Value = '  -2.54%  '

# Row 10
Set-TextValue $ws.Cells.Item(10, 4) '0.06996'
$ws.Cells.Item(10, 5).Value = '  -2.72%  '

# Row 11
Set-TextValue $ws.Cells.Item(11, 4) '0.7783'
$ws.Cells.Item(11, 5).Value = '  -2.55%  '

# Row 12
Set-TextValue $ws.Cells.Item(12, 4) '0.08008'
$ws.Cells.Item(12, 5).Value = '  -0.98%  '

# Row 13
Set-TextValue $ws.Cells.Item(13, 4) '1.932.10'
$ws.Cells.Item(13, 5).Value = '  +0.10%  '

# Row 14
Set-TextValue $ws.Cells.Item(14, 4) '5.344'
$ws.Cells.Item(14, 5).Value = '  -1.47%  '

# Row 15
Set-TextValue $ws.Cells.Item(15, 4) '94.36'
$ws.Cells.Item(15, 5).Value = '  -0.44%  '

# Row 16
Set-TextValue $ws.Cells.Item(16, 4) '14.38'
$ws.Cells.Item(16, 5).Value = '  -3.09%  '

# Row 17
Set-TextValue $ws.Cells.Item(17, 4) '30.324.81'
$ws.Cells.Item(17, 5).Value = '  +0.05%  '

# Row 18
Set-TextValue $ws.Cells.Item(18, 4) '253.07'
$ws.Cells.Item(18, 5).Value = '  +0.68%  '

# Row 19
Set-TextValue $ws.Cells.Item(19, 4) '0.000007922'
$ws.Cells.Item(19, 5).Value = '  -2.64%  '

# Row 20
Set-TextValue $ws.Cells.Item(20, 4) '5.740'
$ws.Cells.Item(20, 5).Value = '  -0.79%  '

# Row 21
Set-TextValue $ws.Cells.Item(21, 4) '2.191.61'
$ws.Cells.Item(21, 5).Value = '  +0.56%  '

# Row 22
Set-TextValue $ws.Cells.Item(22, 4) '0.9991'
$ws.Cells.Item(22, 5).Value = '  -0.21%  '

# Row 23
Set-TextValue $ws.Cells.Item(23, 4) '0.9996'
$ws.Cells.Item(23, 5).Value = '  -0.22%  '

# Row 24
Set-TextValue $ws.Cells.Item(24, 4) '6.660'
$ws.Cells.Item(24, 5).Value = '  -3.60%  '

# Row 25
Set-TextValue $ws.Cells.Item(25, 4) '9.459'
$ws.Cells.Item(25, 5).Value = '  -2.71%  '

# Row 26
Set-TextValue $ws.Cells.Item(26, 4) '165.69'
$ws.Cells.Item(26, 5).Value = '  +0.58%  '

# Row 27
Set-TextValue $ws.Cells.Item(27, 4) '18.99'
$ws.Cells.Item(27, 5).Value = '  -1.12%  '

# Row 28
$ws.Cells.Item(28, 5).Value = '  +3.99%  '

# Row 29
Set-TextValue $ws.Cells.Item(29, 4) '2.192'
$ws.Cells.Item(29, 5).Value = '  -5.31%  '

# Row 30
$ws.Cells.Item(30, 5).Value = '  +0.38%  '

# Row 31
$ws.Cells.Item(31, 5).Value = '  -1.83%  '

# Row 32
Set-TextValue $ws.Cells.Item(32, 4) '4.394'
$ws.Cells.Item(32, 5).Value = '  -0.77%  '

# Row 33
Set-TextValue $ws.Cells.Item(33, 4) '4.120'
$ws.Cells.Item(33, 5).Value = '  -1.99%  '

# Row 34
Set-TextValue $ws.Cells.Item(34, 4) '0.05147'
$ws.Cells.Item(34, 5).Value = '  -1.22%  '

# Row 35
Set-TextValue $ws.Cells.Item(35, 4) '1.284'
$ws.Cells.Item(35, 5).Value = '  +1.33%  '

# Row 36
Set-TextValue $ws.Cells.Item(36, 4) '0.7510'
$ws.Cells.Item(36, 5).Value = '  +0.56%  '

# Row 37
Set-TextValue $ws.Cells.Item(37, 4) '2.770'
$ws.Cells.Item(37, 5).Value = '  +0.35%  '

# Row 38
Set-TextValue $ws.Cells.Item(38, 4) '0.01958'
$ws.Cells.Item(38, 5).Value = '  +0.15%  '

# Row 39
Set-TextValue $ws.Cells.Item(39, 4) '2.800'
$ws.Cells.Item(39, 5).Value = '  +0.08%  '

# Row 40
Set-TextValue $ws.Cells.Item(40, 4) '77.36'
$ws.Cells.Item(40, 5).Value = '  -1.94%  '

# Row 41
Set-TextValue $ws.Cells.Item(41, 4) '6.400'
$ws.Cells.Item(41, 5).Value = '  -0.24%  '

# Row 42
Set-TextValue $ws.Cells.Item(42, 4) '0.4447'
$ws.Cells.Item(42, 5).Value = '  -1.71%  '

# Row 43
Set-TextValue $ws.Cells.Item(43, 4) '1.965'
$ws.Cells.Item(43, 5).Value = '  -2.98%  '

# Row 44
Set-TextValue $ws.Cells.Item(44, 4) '1.000'
$ws.Cells.Item(44, 5).Value = '  -0.11%  '

# Row 45
Set-TextValue $ws.Cells.Item(45, 4) '0.8333'
$ws.Cells.Item(45, 5).Value = '  -0.75%  '

# Row 46
Set-TextValue $ws.Cells.Item(46, 4) '100.78'
$ws.Cells.Item(46, 5).Value = '  -0.92%  '

# Row 47
Set-TextValue $ws.Cells.Item(47, 4) '9.803'
$ws.Cells.Item(47, 5).Value = '  +0.33%  '

# Row 48
$ws.Cells.Item(48, 5).Value = '  +1.05%  '

# Row 49
Set-TextValue $ws.Cells.Item(49, 4) '37.42'
$ws.Cells.Item(49, 5).Value = '  +2.23%  '

# Row 50
Set-TextValue $ws.Cells.Item(50, 4) '980.62'
$ws.Cells.Item(50, 5).Value = '  +11.32%  '

# Row 51
$ws.Cells.Item(51, 2).Value = 'Algorand'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws.Cells.Item(51, 4) '0.1179'
$ws.Cells.Item(51, 5).Value = '  +3.83%  '
